$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 11
$ws.Range("H17").Value = 11

$ws.Range("E18").Value = 92
$ws.Range("F18").Value = 28
$ws.Range("H18").Value = 28
